# Refresh the crypto price/volume snapshot (Price column D, Volume(1h) column E),
# including two rank swaps: rows 28/29 (LidoDAOToken <-> WrappedliquidstakedEther2.0)
# and rows 40/41 (Aptos <-> TheSandbox).
#
# Price-column values that look like plain decimals (e.g. "1.003") are written
# with a leading apostrophe and then Style reset to "Normal" so they land back
# as plain General-formatted text cells (matching the source data, which stores
# every Price/Volume cell as text) instead of being auto-coerced into numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.513.49'
$ws.Range("E2").Value = '  +1.51%  '
$ws.Range("D3").Value = '1.825.80'
$ws.Range("E3").Value = '  +1.44%  '
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = "'315.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.46%  '
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").Value = "'0.5183"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -4.65%  '
$ws.Range("D8").Value = "'0.3938"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +3.56%  '
$ws.Range("D9").Value = "'0.07720"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.39%  '
$ws.Range("D10").Value = "'41.95"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.11%  '
$ws.Range("D11").Value = "'1.113"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.73%  '
$ws.Range("E12").Value = '  +3.01%  '
$ws.Range("D13").Value = "'6.285"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.14%  '
$ws.Range("D14").Value = "'1.003"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = "'7.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("D16").Value = '1.823.17'
$ws.Range("E16").Value = '  +1.59%  '
$ws.Range("D17").Value = "'93.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.87%  '
$ws.Range("E18").Value = '  +1.49%  '
$ws.Range("D19").Value = "'0.06622"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.59%  '
$ws.Range("D20").Value = "'17.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.49%  '
$ws.Range("E21").Value = '  +0.12%  '
$ws.Range("D22").Value = "'6.055"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.99%  '
$ws.Range("D23").Value = '28.526.11'
$ws.Range("E23").Value = '  +1.44%  '
$ws.Range("E24").Value = '  -0.23%  '
$ws.Range("D25").Value = "'2.243"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.47%  '
$ws.Range("D26").Value = "'157.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = '  +1.19%  '
$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").Value = "'2.432"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.82%  '
$ws.Range("B29").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C29").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D29").Value = '2.034.26'
$ws.Range("E29").Value = '  +1.54%  '
$ws.Range("D30").Value = "'125.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.53%  '
$ws.Range("D31").Value = "'1.132"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").Value = "'0.1106"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.07%  '
$ws.Range("D33").Value = "'5.666"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.17%  '
$ws.Range("E34").Value = '  -0.55%  '
$ws.Range("D35").Value = "'0.07227"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.40%  '
$ws.Range("D36").Value = "'0.2240"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.40%  '
$ws.Range("D37").Value = "'8.944"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.16%  '
$ws.Range("D38").Value = "'0.02333"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.61%  '
$ws.Range("D39").Value = "'5.165"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.62%  '
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = "'0.6249"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("B41").Value = 'Aptos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D41").Value = "'11.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.85%  '
$ws.Range("E42").Value = '  +1.36%  '
$ws.Range("E43").Value = '  +0.13%  '
$ws.Range("D44").Value = "'1.395"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.67%  '
$ws.Range("E45").Value = '  +0.36%  '
$ws.Range("D46").Value = "'0.5913"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.05%  '
$ws.Range("D47").Value = "'3.704"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.51%  '
$ws.Range("D48").Value = "'124.36"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.03%  '
$ws.Range("D49").Value = "'1.981"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.31%  '
$ws.Range("D50").Value = "'1.186"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.16%  '
$ws.Range("D51").Value = "'0.06941"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.86%  '
